$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'244.50"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'5.403"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'0.06038"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'0.8140"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.9192"
$c.Style = "Normal"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c = $ws.Range("D9")
$c.Value = "'0.1436"
$c.Style = "Normal"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c = $ws.Range("D10")
$c.Value = "'0.07484"
$c.Style = "Normal"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$c = $ws.Range("D11")
$c.Value = "'0.03394"
$c.Style = "Normal"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c = $ws.Range("D12")
$c.Value = "'0.03046"
$c.Style = "Normal"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c = $ws.Range("D13")
$c.Value = "'0.09423"
$c.Style = "Normal"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$c = $ws.Range("D14")
$c.Value = "'4.012"
$c.Style = "Normal"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c = $ws.Range("D15")
$c.Value = "'0.001588"
$c.Style = "Normal"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$c = $ws.Range("D16")
$c.Value = "'0.04827"
$c.Style = "Normal"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$c = $ws.Range("D17")
$c.Value = "'0.0005942"
$c.Style = "Normal"
$ws.Range("E17").Value = "16OneONE"
$c = $ws.Range("D18")
$c.Value = "'0.005697"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'0.004167"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'0.0009901"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'3.661"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'6.436"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'2.183"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'0.1322"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'0.00008409"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.03991"
$c.Style = "Normal"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$c = $ws.Range("D41")
$c.Value = "'0.1076"
$c.Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$c = $ws.Range("D42")
$c.Value = "'0.002723"
$c.Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$c = $ws.Range("D43")
$c.Value = "'0.003051"
$c.Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICK"
$c = $ws.Range("D44")
$c.Value = "'0.005800"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.00005227"
$c.Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$c = $ws.Range("D48")
$c.Value = "'0.002319"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.00002100"
$c.Style = "Normal"
